$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the cell values first
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the full target style (bold font, thin box border, centered/top
# aligned) on B1 in one pass so only a single new cellXf is created.
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108
$b1.VerticalAlignment = -4160
$b1.Borders.LineStyle = 1

# Copy that finished format onto A2 via PasteSpecial (formats only) so A2
# reuses the same style index instead of walking through the same chain
# of intermediate styles again (which would leave orphan cellXfs behind).
$b1.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
